# Renumber column A (index column) in the active sheet so that instead of
# resetting to 0 at the start of each month, it becomes one continuous
# running count for the whole year, starting at 0 on row 2.
#
# Rows 2-23 (January) already hold the correct running values (0..21), so
# only rows 24-248 (Feb..Dec) actually change, but we recompute the whole
# column for safety/consistency.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 248

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 1).Value = $row - $firstRow
}
